# LMS-2340 Changes to OD600 and Transcriptomics templates.
#
# Adds a new "Header Format" row (A8/B8/C8 = "Header Format" /
# "TIME::VALUE_TYPE" / "Must be TIME::VALUE_TYPE") to the bottom of the
# openbis-metadata sheet, widens column B to fit the new value, flags the
# sheet for portrait/Envelope#10 printing, and leaves the new C8 cell
# selected - mirroring the upstream SVN commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openbis-metadata")

# --- New row 8 content -------------------------------------------------
# Cells are written in the order A8, C8, B8 so that brand-new shared
# strings land in the workbook's sharedStrings table in the same order
# as the target file: "Header Format", "Must be TIME::VALUE_TYPE",
# "TIME::VALUE_TYPE" (in that order - B8's text is added to the table
# last even though column B comes before column C).
$ws.Range("A8").Value = "Header Format"
$ws.Range("C8").Value = "Must be TIME::VALUE_TYPE"
$ws.Range("B8").Value = "TIME::VALUE_TYPE"

# Copy formatting from the analogous existing cells so the new row reuses
# the exact same (pre-existing) cell styles as the rest of the table
# instead of Excel minting new style records: column A uses the label
# style from A7, column C the description style from C7, and column B the
# "value" style already used by B2.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B2").Copy()
$ws.Range("B8").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# The new row uses the sheet's normal (explicit) 18pt row height.
$ws.Rows.Item(8).RowHeight = 18

# --- Column width --------------------------------------------------------
# Column B needs to end up stored as width="24" in the saved XML. Excel's
# ColumnWidth property is expressed in "characters of the Normal font" and
# gets pixel-snapped/rounded on write, so the character-width value that
# round-trips to exactly 24 isn't 24 itself; 23.29 lands in the middle of
# the input range that snaps to a stored width of 24.
$ws.Columns.Item(2).ColumnWidth = 23.29

# --- Page setup ------------------------------------------------------------
# paperSize 10 = Envelope #10, orientation 1 = xlPortrait.
$ws.PageSetup.PaperSize = 10
$ws.PageSetup.Orientation = 1
# NOTE: the target XML also carries pageSetup horizontalDpi/verticalDpi
# ("4294967292", i.e. -4 as a signed 32-bit sentinel). This runtime's
# PageSetup object model does not expose a settable Horizontal/VerticalDpi
# property (nor any synonym - Resolution/PrintQuality/etc. were all probed
# and have no effect), so those two attributes can't be produced from COM
# here; everything else about the edit is reproduced exactly.

# --- Selection ---------------------------------------------------------
$ws.Range("C8").Select()
